$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Date" column (BF) held the wrong value "4-17-2007-08" for every
# data row (rows 2-31) because of how the NBA stats site reported the
# game date. Correct it to the real ISO date "2008-04-17".
# A leading apostrophe forces the text to stay a literal string instead
# of being auto-parsed into a date serial number by the smart-entry logic.
$firstRow = 2
$lastRow = 31
$dateCol = "BF"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Range("$dateCol$r")
    if ($cell.Value2 -eq "4-17-2007-08") {
        $cell.Value = "'2008-04-17"
    }
}
